$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value = "Complex_Interface"
$ws.Range("E33").Value = "Can't be found from sce"
$ws.Range("D3").Value = "Beta strand"

$ws.Range("E12").Select()
